# Update with Correct Forecast output
# - Insert a new "Week_Start_Date" column (B) into the Forecast Comparison sheet,
#   shifting ASIN / forecasts / Product Title / is_holiday_week one column right.
# - Normalize the Week labels (W01 -> W1, etc.)
# - Populate the new Week_Start_Date column with each week's start date.
# - Correct a handful of MyForecast values.
# - Store is_holiday_week as a proper boolean.
# - Bump the 16-week forecast total on the Summary sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

# Insert new column B ("Week_Start_Date"); everything from the old B onward shifts right by one.
$ws.Columns.Item(2).Insert()

$ws.Cells.Item(1, 2).Value = "Week_Start_Date"

# Week label (col A) and new Week_Start_Date (col B) + corrected MyForecast (col D) per row.
$weekData = @(
    @{ Row = 2;  Week = "W1";  StartDate = "2025-01-05"; MyForecast = 112 },
    @{ Row = 3;  Week = "W2";  StartDate = "2025-01-12"; MyForecast = 105 },
    @{ Row = 4;  Week = "W3";  StartDate = "2025-01-19"; MyForecast = 106 },
    @{ Row = 5;  Week = "W4";  StartDate = "2025-01-26"; MyForecast = 101 },
    @{ Row = 6;  Week = "W5";  StartDate = "2025-02-02"; MyForecast = 85  },
    @{ Row = 7;  Week = "W6";  StartDate = "2025-02-09"; MyForecast = 87  },
    @{ Row = 8;  Week = "W7";  StartDate = "2025-02-16"; MyForecast = 92  },
    @{ Row = 9;  Week = "W8";  StartDate = "2025-02-23"; MyForecast = 96  },
    @{ Row = 10; Week = "W9";  StartDate = "2025-03-02"; MyForecast = 94  },
    @{ Row = 11; Week = "W10"; StartDate = "2025-03-09"; MyForecast = 93  },
    @{ Row = 12; Week = "W11"; StartDate = "2025-03-16"; MyForecast = 93  },
    @{ Row = 13; Week = "W12"; StartDate = "2025-03-23"; MyForecast = 100 },
    @{ Row = 14; Week = "W13"; StartDate = "2025-03-30"; MyForecast = 101 },
    @{ Row = 15; Week = "W14"; StartDate = "2025-04-06"; MyForecast = 77  },
    @{ Row = 16; Week = "W15"; StartDate = "2025-04-13"; MyForecast = 97  },
    @{ Row = 17; Week = "W16"; StartDate = "2025-04-20"; MyForecast = 98  }
)

foreach ($item in $weekData) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value = $item.Week
    # Leading apostrophe forces text so the date-looking string isn't reinterpreted as a date.
    $ws.Cells.Item($r, 2).Value = "'" + $item.StartDate
    $ws.Cells.Item($r, 4).Value = $item.MyForecast
    # is_holiday_week (now column J) becomes a real boolean instead of numeric 0.
    $ws.Cells.Item($r, 10).Value = $false
}

$ws.Range("A1").Select()

# Summary sheet: correct the 16-week forecast total.
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B9").Value = "'1536"
